# Merge the split "<id>" / "p038v_N" / "</id>" runs back into a single run
# for each of the six <id> elements in the document (p038v_1 .. p038v_6).
#
# Before: three runs -> <id>  (Courier formatting)
#                        p038v_N  (plain black formatting)
#                        </id>  (Courier formatting)
# After:  one run     -> <id>p038v_N</id>  (Courier formatting, same as the
#                          original surrounding "<id>" / "</id>" runs)
#
# A simple Find/Replace over the concatenated text merges the three runs
# into one, picking up the formatting of the first run in the matched
# range (the Courier New "<id>" run), which matches the target formatting
# exactly.

$d = $word.ActiveDocument

for ($i = 1; $i -le 6; $i++) {
    $tag = "<id>p038v_$i</id>"
    $d.Content.Find.Execute($tag, $true, $false, $false, $false, $false, $true, 1, $false, $tag, 2) | Out-Null
}
